$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.702.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "'3.193.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.01%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'589.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("D6").Value = "'136.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'3.188.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.08%  "
$ws.Range("D9").Value = "'0.508"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("E10").Value = "  -4.38%  "
$ws.Range("D11").Value = "'5.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.42%  "
$ws.Range("D12").Value = "'0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.32%  "
$ws.Range("D13").Value = "'0.0000237"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.09%  "
$ws.Range("D14").Value = "'33.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("D15").Value = "'3.712.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.22%  "
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").Value = "'3.187.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.24%  "
$ws.Range("D18").Value = "'62.695.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").Value = "'6.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.36%  "
$ws.Range("D20").Value = "'457.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.54%  "
$ws.Range("D21").Value = "'13.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value = "'0.706"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.41%  "
$ws.Range("D23").Value = "'7.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.19%  "
$ws.Range("D24").Value = "'13.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'83.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("E29").Value = "  -6.04%  "
$ws.Range("D30").Value = "'7.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.80%  "
$ws.Range("E31").Value = "  -6.16%  "
$ws.Range("D32").Value = "'27.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.55%  "
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").Value = "'2.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("E35").Value = "  -5.38%  "
$ws.Range("E36").Value = "  -0.49%  "
$ws.Range("D37").Value = "'51.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.20%  "
$ws.Range("D38").Value = "'0.0₃0701"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.15%  "
$ws.Range("D39").Value = "'0.0387"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'2.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'402.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "'8.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.23%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'2.852.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.27%  "
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").Value = "'36.89"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.72%  "
$ws.Range("E46").Value = "  -5.57%  "
$ws.Range("D47").Value = "'2.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "'125.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").Value = "'25.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("E51").Value = "  -3.11%  "
